$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($rangeAddress, $text) {
    # Binance's OHLCV numbers are stored as literal text (not numbers) in
    # this sheet. Force the Text number format before writing so a
    # numeric-looking string ("142.46000000") isn't auto-coerced into a
    # float, then drop back to the default cell style so we don't leave a
    # visible NumberFormat behind (matches how the rest of the sheet looks).
    $rng = $ws.Range($rangeAddress)
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.Style = "Normal"
}

# --- 1. Open_time (column B) was stored in seconds; convert to milliseconds
#        for rows 2-16 (multiply existing value by 1000). ---
for ($r = 2; $r -le 16; $r++) {
    $cell = $ws.Cells.Item($r, 2)
    $seconds = $cell.Value()
    $cell.Value = $seconds * 1000
}

# --- 2. Row 16 got re-computed with the (now complete) daily candle data. ---
Set-TextValue "F16" "142.46000000"
Set-TextValue "G16" "151189.65877000"
Set-TextValue "I16" "21495196.81888500"
$ws.Range("J16").Value = 77131
Set-TextValue "K16" "73238.69048000"
Set-TextValue "L16" "10413387.06460380"

# --- 3. A new row (17) was appended for the next day's candle. Copy row 16
#        first so the new row inherits the same cell styles (bold+border
#        style on column A), then overwrite every value with the real data
#        for the new row. ---
$ws.Range("A16:M16").Copy($ws.Range("A17:M17"))

$ws.Range("A17").Value = 15
$ws.Range("B17").Value = 1576454400000
Set-TextValue "C17" "142.46000000"
Set-TextValue "D17" "142.72000000"
Set-TextValue "E17" "140.30000000"
Set-TextValue "F17" "141.12000000"
Set-TextValue "G17" "78497.62898000"
$ws.Range("H17").Value = 1576540799999
Set-TextValue "I17" "11089010.61664710"
$ws.Range("J17").Value = 41602
Set-TextValue "K17" "33864.95452000"
Set-TextValue "L17" "4785822.81632040"
Set-TextValue "M17" "0"

Write-Output "done"
